# Apply the phenotype header re-casing and updated "endometritis" (column C)
# probability values as described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (row 1) ---
$ws.Range("B1").Value = "Pelvic Inflammatory Diseases"
$ws.Range("C1").Value = "Endometritis"
$ws.Range("D1").Value = "Miscarriage"
$ws.Range("E1").Value = "Ovarian Cancer"
$ws.Range("F1").Value = "Cervical Intraepithelial Neoplasia"
$ws.Range("G1").Value = "Preterm Prelabor Rupture of Membranes (PPROM)"
$ws.Range("H1").Value = "Dysmenorrhea (Menstrual pain)"
$ws.Range("I1").Value = "Adenomyosis"
$ws.Range("J1").Value = "Vaginal Dryness"
$ws.Range("K1").Value = "Gestational Diabetes"

# --- Updated "Endometritis" (column C) score values ---
$ws.Range("C2").Value = 37.4
$ws.Range("C3").Value = 6.7
$ws.Range("C4").Value = 61.5
$ws.Range("C5").Value = 66.3
$ws.Range("C7").Value = 10.4
$ws.Range("C8").Value = 5.6
$ws.Range("C9").Value = 23.4
$ws.Range("C11").Value = 5.9
$ws.Range("C12").Value = 42.7
$ws.Range("C13").Value = 38.9
$ws.Range("C14").Value = 74.90000000000001
$ws.Range("C15").Value = 80.3
$ws.Range("C16").Value = 60.7
$ws.Range("C17").Value = 72.2
$ws.Range("C18").Value = 83.09999999999999
$ws.Range("C19").Value = 60
$ws.Range("C20").Value = 69.40000000000001
$ws.Range("C21").Value = 95
$ws.Range("C22").Value = 87.2
$ws.Range("C23").Value = 89.8
$ws.Range("C24").Value = 44.9
$ws.Range("C25").Value = 19.9
$ws.Range("C26").Value = 68.90000000000001
$ws.Range("C27").Value = 6.4
$ws.Range("C28").Value = 30.8
$ws.Range("C29").Value = 54.2
$ws.Range("C30").Value = 74.5
$ws.Range("C33").Value = 31
$ws.Range("C34").Value = 42.1
$ws.Range("C35").Value = 90.09999999999999
$ws.Range("C37").Value = 69.40000000000001
$ws.Range("C38").Value = 47.2
$ws.Range("C39").Value = 65.09999999999999
$ws.Range("C40").Value = 62.2
$ws.Range("C41").Value = 25.5
$ws.Range("C42").Value = 7.9
$ws.Range("C43").Value = 25.2
$ws.Range("C44").Value = 30.5
$ws.Range("C45").Value = 83
$ws.Range("C46").Value = 33.1
$ws.Range("C47").Value = 54.1
$ws.Range("C48").Value = 94.09999999999999
$ws.Range("C49").Value = 42.7
$ws.Range("C50").Value = 70.3
$ws.Range("C51").Value = 29.2
$ws.Range("C52").Value = 22.9
$ws.Range("C53").Value = 69.40000000000001
$ws.Range("C54").Value = 82.09999999999999
$ws.Range("C55").Value = 20.2
$ws.Range("C56").Value = 31.3
$ws.Range("C57").Value = 50.8
